# "commit to merge files" - append the merged firstname/lastname/email/
# confirm-pass/mobileno credential row (with its header row) below the
# existing username/password table on the Credentials sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

# Populate the new header (row 6) and data (row 7) cells in the same
# order the strings first appear in the authored workbook, so the
# resulting shared-strings table comes out in the same sequence.
$ws.Range("D7").Value = "prati23"
$ws.Range("B6").Value = "firstname"
$ws.Range("B7").Value = "pratibha"
$ws.Range("C6").Value = "lastname"
$ws.Range("C7").Value = "wabale"
$ws.Range("D6").Value = "email"
$ws.Range("F6").Value = "confirm pass"
$ws.Range("H6").Value = "mobileno"
$ws.Range("E7").Value = "EnterPassword@1234"
$ws.Range("F7").Value = "EnterPassword@1234"

# These reuse shared strings already present in the workbook.
$ws.Range("A6").Value = "username"
$ws.Range("E6").Value = "password"
$ws.Range("A7").Value = "pratipw"
$ws.Range("H7").Value = 6754902280

# Resize the columns to fit the newly-added content (closest width the
# host's character->width conversion can reach to the authored values).
$ws.Columns.Item(3).ColumnWidth = 13
$ws.Columns.Item(4).ColumnWidth = 22.666666666666668
$ws.Columns.Item(6).ColumnWidth = 18.833333333333332
$ws.Columns.Item(8).ColumnWidth = 10

# Scroll the view over and leave the new mobile-number cell selected,
# matching the author's final on-screen state.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H7").Select() | Out-Null
